$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 20
$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 30

$ws.Range("C9").Select()
